# UC1_TC2 evaluation workbook update — refresh the QuantitativeMetrics sheet
# with the latest run results (runtime exception found, CodeBLEU re-scored).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")
$ws.Activate()

# --- Execution metrics --------------------------------------------------

# Runtime without error: now fails, with an explanatory note.
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "Runtime exception due to invalid locator"

# Assertion validity: previously "no" with a note; now cleared/blank
# because the run never reached the assertion stage.
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = ""

# --- Syntax & Semantic similarity ---------------------------------------

# Updated CodeBLEU score + breakdown (dataflow_match_score improved).
$ws.Range("B12").Value = 0.2140601269782228
$ws.Range("C12").Value = "{'codebleu': 0.21406012697822277, 'ngram_match_score': 0.0753538535441081, 'weighted_ngram_match_score': 0.10968278262044664, 'syntax_match_score': 0.5045372050816697, 'dataflow_match_score': 0.16666666666666666}"

# --- Selection / cursor position -----------------------------------------
$ws.Range("C8").Select()
